# LogBook.xlsx update:
#  - Row 15 "Finish Date" (B15) corrected from 30/04/2021 to 27/04/2021.
#  - A new log entry added on row 16 (28/04/2021 - 28/04/2021) describing
#    final report / table / comment work, reusing the same date style as
#    the row above it (so no new number-format style is created) and the
#    existing wrap-text style already present on C16/D16.
#  - Selection moved to the newly entered D16 cell, matching where the
#    user ended up after typing the new row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- correct the Finish Date on the existing last entry (row 15) ---
$ws.Range("B15").Value = 44313

# --- bring in the date formatting from row 15 so A16/B16 share style s="2" ---
$ws.Range("A15:B15").Copy() | Out-Null
$ws.Range("A16:B16").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws.Range("A16").Value = 44314
$ws.Range("B16").Value = 44314
$ws.Range("C16").Value = "I need to improve my report for final submission. I need to add tables. I need to review the code and make changes if needed."
$ws.Range("D16").Value = "Tables have been added to the report.  Comments have been added to the code to help the user understand the code better."

# row grew to fit the wrapped notes, same as the other multi-line rows
$ws.Rows.Item(16).RowHeight = 43.5

# leave the cursor where the author left it after finishing the new row
$ws.Range("D16").Select() | Out-Null
